$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Combined Calibration")

# New "VSA Offset (dB)" values for E2:E82 (column E), recomputed per
# the 20250929 refactor pass on the combined calibration data.
$evmValues = @{
    2 = 1.6384837962300001
    3 = 1.6586218803099999
    4 = 1.69783757327
    5 = 1.66663520516
    6 = 1.6394822981399999
    7 = 1.6707124931300001
    8 = 1.5562517593499998
    9 = 1.6695577590899999
    10 = 1.6901819344600002
    11 = 1.6327185529300001
    12 = 1.6529666815700002
    13 = 1.6806430130099999
    14 = 1.7207988732899999
    15 = 1.7680798639299999
    16 = 1.8414640190599998
    17 = 1.84070112484
    18 = 1.8190075346300001
    19 = 1.78324713921
    20 = 1.8669747491299999
    21 = 1.8833715236700002
    22 = 1.7462585513500002
    23 = 1.8775014061699999
    24 = 1.8444125277249999
    25 = 1.82673410661
    26 = 1.9237036086000001
    27 = 1.8607387499650001
    28 = 1.83742717552
    29 = 1.96477983914
    30 = 1.9203338963600001
    31 = 1.9600063868099999
    32 = 2.1004507234199998
    33 = 2.0897662510199999
    34 = 1.9686865043699999
    35 = 1.98511326897
    36 = 2.0487600563999999
    37 = 2.0188485269599998
    38 = 2.0895582085600002
    39 = 2.1455872832999998
    40 = 2.0128021831899998
    41 = 1.9354800532399998
    42 = 2.01036161066
    43 = 2.1985707115999999
    44 = 2.0687919849500003
    45 = 2.1388264352799999
    46 = 2.19407232391
    47 = 2.1778569245499999
    48 = 2.2173617492000002
    49 = 2.2300445945599998
    50 = 2.1968254528358999
    51 = 2.24478663353
    52 = 2.2341556567500001
    53 = 2.2113324208139997
    54 = 2.1556990335499999
    55 = 2.2393885515200003
    56 = 2.2712320039099998
    57 = 2.2958746626240001
    58 = 2.2380031172500003
    59 = 2.1561174938400001
    60 = 2.3908738164700001
    61 = 2.3444008994765997
    62 = 2.32682532251
    63 = 2.4198693923870001
    64 = 2.7707788727999998
    65 = 2.6148606409999999
    66 = 2.5843345257000001
    67 = 2.5671409818999997
    68 = 2.4475936934
    69 = 2.5853386400999998
    70 = 2.5519596187999998
    71 = 2.5514894947000002
    72 = 2.5813201328000002
    73 = 2.6591293584000004
    74 = 2.6963444582
    75 = 2.5840228168999997
    76 = 2.6822460193
    77 = 2.6878703581000001
    78 = 2.7332143242
    79 = 2.7355281550999999
    80 = 2.7539368351999998
    81 = 2.7657501017000001
    82 = 2.7365483840000002
}

foreach ($row in $evmValues.Keys) {
    $ws.Cells.Item($row, 5).Value = $evmValues[$row]
}

